$d = $word.ActiveDocument

# 1) The "Wiesent" bridge description was split across three runs so that a
#    spell-check proofErr pair could bracket the word "Wiesent". Collapse it
#    back down to a single plain run (no proofErr, no mid-sentence run
#    breaks) by replacing the whole sentence with itself.
$d.Content.Find.Execute(
    "Start in Wadendorf an der Brücke über die Wiesent in Richtung Scherleithen (westlich, leicht ansteigend)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Start in Wadendorf an der Brücke über die Wiesent in Richtung Scherleithen (westlich, leicht ansteigend)",
    2) | Out-Null

# 2) Append a new paragraph (after the existing trailing empty paragraph)
#    with the GPS coordinates of the turning point, split across four runs.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$gpsPara = $d.Paragraphs.Last
$gpsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">GPS (Wendeplatz): N </w:t></w:r>' +
    '<w:r><w:t>49.8895309,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> E </w:t></w:r>' +
    '<w:r><w:t>11.3135538</w:t></w:r>' +
    '</w:p>'
$gpsPara.Range.InsertXML($gpsXml) | Out-Null
